$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.052410929767802
$ws.Cells.Item(2, 3).Value = 0.1975470831521875
$ws.Cells.Item(2, 4).Value = 0.02126524561209209
$ws.Cells.Item(2, 5).Value = 0.08822327336754476
$ws.Cells.Item(2, 6).Value = 0.6878680425175787
$ws.Cells.Item(2, 9).Value = 0.5845085108126327
$ws.Cells.Item(2, 12).Value = 0.220872605439709
$ws.Cells.Item(2, 14).Value = 1.079996518808514
$ws.Cells.Item(2, 15).Value = 2.325126003310601
$ws.Cells.Item(3, 2).Value = 0.9510981902345179
$ws.Cells.Item(3, 3).Value = 0.1841607495851179
$ws.Cells.Item(3, 4).Value = 0.01995137572605898
$ws.Cells.Item(3, 5).Value = 0.08879252124508064
$ws.Cells.Item(3, 6).Value = 0.6814553148613385
$ws.Cells.Item(3, 9).Value = 0.5892430574794965
$ws.Cells.Item(3, 12).Value = 0.211234034076071
$ws.Cells.Item(3, 14).Value = 1.087544628982975
$ws.Cells.Item(3, 15).Value = 2.317957044106976
$ws.Cells.Item(4, 2).Value = 0.8889802177097863
$ws.Cells.Item(4, 3).Value = 0.1758770014104414
$ws.Cells.Item(4, 4).Value = 0.01913838266573009
$ws.Cells.Item(4, 5).Value = 0.08918658682360459
$ws.Cells.Item(4, 6).Value = 0.6780108799843703
$ws.Cells.Item(4, 9).Value = 0.5925344515371762
$ws.Cells.Item(4, 12).Value = 0.2054272173411107
$ws.Cells.Item(4, 14).Value = 1.092619227828692
$ws.Cells.Item(4, 15).Value = 2.315197429004968
$ws.Cells.Item(5, 2).Value = 0.863690449485091
$ws.Cells.Item(5, 3).Value = 0.1724852679513589
$ws.Cells.Item(5, 4).Value = 0.01880552398137425
$ws.Cells.Item(5, 5).Value = 0.08935837718386708
$ws.Cells.Item(5, 6).Value = 0.6767311016908693
$ws.Cells.Item(5, 9).Value = 0.5939722641412644
$ws.Cells.Item(5, 12).Value = 0.203088963409698
$ws.Cells.Item(5, 14).Value = 1.094797927558623
$ws.Cells.Item(5, 15).Value = 2.314485336073773
$ws.Cells.Item(6, 2).Value = 0.8594925869730048
$ws.Cells.Item(6, 3).Value = 0.1719211096941535
$ws.Cells.Item(6, 4).Value = 0.01875015956515469
$ws.Cells.Item(6, 5).Value = 0.08938757974779676
$ws.Cells.Item(6, 6).Value = 0.6765260728177367
$ws.Cells.Item(6, 9).Value = 0.5942168396013123
$ws.Cells.Item(6, 12).Value = 0.2027023959845735
$ws.Cells.Item(6, 14).Value = 1.095166391931876
$ws.Cells.Item(6, 15).Value = 2.314391992456308
$ws.Cells.Item(7, 2).Value = 0.8886390525689194
$ws.Cells.Item(7, 3).Value = 0.1758313240119094
$ws.Cells.Item(7, 4).Value = 0.01913389989611858
$ws.Cells.Item(7, 5).Value = 0.08918885827237411
$ws.Cells.Item(7, 6).Value = 0.6779931191065103
$ws.Cells.Item(7, 9).Value = 0.5925534516275022
$ws.Cells.Item(7, 12).Value = 0.2053955690619063
$ws.Cells.Item(7, 14).Value = 1.092648161925439
$ws.Cells.Item(7, 15).Value = 2.315186155999982
$ws.Cells.Item(8, 2).Value = 1.017460844502466
$ws.Cells.Item(8, 3).Value = 0.1929449612804888
$ws.Cells.Item(8, 4).Value = 0.02081353636032901
$ws.Cells.Item(8, 5).Value = 0.08841030716870435
$ws.Cells.Item(8, 6).Value = 0.6855545383949817
$ws.Cells.Item(8, 9).Value = 0.5860611422964368
$ws.Cells.Item(8, 12).Value = 0.2175261706784397
$ws.Cells.Item(8, 14).Value = 1.082507844134042
$ws.Cells.Item(8, 15).Value = 2.322313030469701
$ws.Cells.Item(9, 2).Value = 1.270725647457823
$ws.Cells.Item(9, 3).Value = 0.2259868145100086
$ws.Cells.Item(9, 4).Value = 0.02405681686477124
$ws.Cells.Item(9, 5).Value = 0.08723686863369551
$ws.Cells.Item(9, 6).Value = 0.7043012336873318
$ws.Cells.Item(9, 9).Value = 0.5763846500694285
$ws.Cells.Item(9, 12).Value = 0.2421952227594772
$ws.Cells.Item(9, 14).Value = 1.066109642971981
$ws.Cells.Item(9, 15).Value = 2.349343833708645
$ws.Cells.Item(10, 2).Value = 1.45713589541549
$ws.Cells.Item(10, 3).Value = 0.2499406610461392
$ws.Cells.Item(10, 4).Value = 0.0264081263776248
$ws.Cells.Item(10, 5).Value = 0.08658998129950035
$ws.Cells.Item(10, 6).Value = 0.7204758102281659
$ws.Cells.Item(10, 9).Value = 0.5711452031298094
$ws.Cells.Item(10, 12).Value = 0.2608558580581644
$ws.Cells.Item(10, 14).Value = 1.056182097808296
$ws.Cells.Item(10, 15).Value = 2.377203493937657
$ws.Cells.Item(11, 2).Value = 1.542000895782223
$ws.Cells.Item(11, 3).Value = 0.2607667725354474
$ws.Cells.Item(11, 4).Value = 0.02747080463101526
$ws.Cells.Item(11, 5).Value = 0.08634240727371179
$ws.Cells.Item(11, 6).Value = 0.7283582978814138
$ws.Cells.Item(11, 9).Value = 0.5691693153636308
$ws.Cells.Item(11, 12).Value = 0.2694614766756871
$ws.Cells.Item(11, 14).Value = 1.052125141489121
$ws.Cells.Item(11, 15).Value = 2.391624023469006
$ws.Cells.Item(12, 2).Value = 1.574145223253652
$ws.Cells.Item(12, 3).Value = 0.2648560342221344
$ws.Cells.Item(12, 4).Value = 0.02787219736050162
$ws.Cells.Item(12, 5).Value = 0.08625536996305172
$ws.Cells.Item(12, 6).Value = 0.7314187963054053
$ws.Cells.Item(12, 9).Value = 0.5684798350986071
$ws.Cells.Item(12, 12).Value = 0.2727369443705214
$ws.Cells.Item(12, 14).Value = 1.050654824032975
$ws.Cells.Item(12, 15).Value = 2.397336545232321
$ws.Cells.Item(13, 2).Value = 1.567222047657935
$ws.Cells.Item(13, 3).Value = 0.263975802093313
$ws.Cells.Item(13, 4).Value = 0.02778579598800235
$ws.Cells.Item(13, 5).Value = 0.08627381644093646
$ws.Cells.Item(13, 6).Value = 0.7307563000046287
$ws.Cells.Item(13, 9).Value = 0.5686257120839926
$ws.Cells.Item(13, 12).Value = 0.2720307721875912
$ws.Cells.Item(13, 14).Value = 1.05096855033387
$ws.Cells.Item(13, 15).Value = 2.396095044594603
$ws.Cells.Item(14, 2).Value = 1.544645282498038
$ws.Cells.Item(14, 3).Value = 0.2611034071505856
$ws.Cells.Item(14, 4).Value = 0.02750384802961037
$ws.Cells.Item(14, 5).Value = 0.08633511212028999
$ws.Cells.Item(14, 6).Value = 0.7286085715923747
$ws.Cells.Item(14, 9).Value = 0.5691114132668531
$ws.Cells.Item(14, 12).Value = 0.2697306170601763
$ws.Cells.Item(14, 14).Value = 1.052002855894685
$ws.Cells.Item(14, 15).Value = 2.392088946829489
$ws.Cells.Item(15, 2).Value = 1.530817326412887
$ws.Cells.Item(15, 3).Value = 0.2593426283442284
$ws.Cells.Item(15, 4).Value = 0.02733101322064613
$ws.Cells.Item(15, 5).Value = 0.08637353175626217
$ws.Cells.Item(15, 6).Value = 0.7273028715599139
$ws.Cells.Item(15, 9).Value = 0.5694165737751078
$ws.Cells.Item(15, 12).Value = 0.2683238786775064
$ws.Cells.Item(15, 14).Value = 1.052644986684726
$ws.Cells.Item(15, 15).Value = 2.389667902034972
$ws.Cells.Item(16, 2).Value = 1.451590951338915
$ws.Cells.Item(16, 3).Value = 0.2492317113716922
$ws.Cells.Item(16, 4).Value = 0.02633853624700322
$ws.Cells.Item(16, 5).Value = 0.08660710040017783
$ws.Cells.Item(16, 6).Value = 0.7199712388882062
$ws.Cells.Item(16, 9).Value = 0.5712825450968708
$ws.Cells.Item(16, 12).Value = 0.2602958033715339
$ws.Cells.Item(16, 14).Value = 1.05645646323692
$ws.Cells.Item(16, 15).Value = 2.376296281024224
$ws.Cells.Item(17, 2).Value = 1.403003838688335
$ws.Cells.Item(17, 3).Value = 0.2430107647477087
$ws.Cells.Item(17, 4).Value = 0.02572788959970751
$ws.Cells.Item(17, 5).Value = 0.08676234664163474
$ws.Cells.Item(17, 6).Value = 0.7156079773055808
$ws.Cells.Item(17, 9).Value = 0.5725317433550074
$ws.Cells.Item(17, 12).Value = 0.2554006915600411
$ws.Cells.Item(17, 14).Value = 1.058912232329241
$ws.Cells.Item(17, 15).Value = 2.368541080038
$ws.Cells.Item(18, 2).Value = 1.375064134507568
$ws.Cells.Item(18, 3).Value = 0.239426006107351
$ws.Cells.Item(18, 4).Value = 0.02537600948146945
$ws.Cells.Item(18, 5).Value = 0.08685603572949674
$ws.Cells.Item(18, 6).Value = 0.7131477136699829
$ws.Cells.Item(18, 9).Value = 0.5732886026417034
$ws.Cells.Item(18, 12).Value = 0.2525961520362756
$ws.Cells.Item(18, 14).Value = 1.060367944823199
$ws.Cells.Item(18, 15).Value = 2.364244899060054
$ws.Cells.Item(19, 2).Value = 1.365605367813487
$ws.Cells.Item(19, 3).Value = 0.2382111349433558
$ws.Cells.Item(19, 4).Value = 0.02525675755007484
$ws.Cells.Item(19, 5).Value = 0.08688851223300453
$ws.Cells.Item(19, 6).Value = 0.7123231856937338
$ws.Cells.Item(19, 9).Value = 0.5735514452166868
$ws.Cells.Item(19, 12).Value = 0.2516484750015593
$ws.Cells.Item(19, 14).Value = 1.060868248855165
$ws.Cells.Item(19, 15).Value = 2.362818505716433
$ws.Cells.Item(20, 2).Value = 1.408175378240628
$ws.Cells.Item(20, 3).Value = 0.2436736832330268
$ws.Cells.Item(20, 4).Value = 0.02579296164895339
$ws.Cells.Item(20, 5).Value = 0.08674536550199718
$ws.Cells.Item(20, 6).Value = 0.7160673434521385
$ws.Cells.Item(20, 9).Value = 0.5723947936203686
$ws.Cells.Item(20, 12).Value = 0.2559206468274766
$ws.Cells.Item(20, 14).Value = 1.058646338947128
$ws.Cells.Item(20, 15).Value = 2.369349615002449
$ws.Cells.Item(21, 2).Value = 1.551276427734933
$ws.Cells.Item(21, 3).Value = 0.2619473816907032
$ws.Cells.Item(21, 4).Value = 0.02758669084867194
$ws.Cells.Item(21, 5).Value = 0.08631692590063089
$ws.Cells.Item(21, 6).Value = 0.7292373594273585
$ws.Cells.Item(21, 9).Value = 0.5689671556162779
$ws.Cells.Item(21, 12).Value = 0.2704057756758544
$ws.Cells.Item(21, 14).Value = 1.051697265486254
$ws.Cells.Item(21, 15).Value = 2.393258797168016
$ws.Cells.Item(22, 2).Value = 1.64484587520019
$ws.Cells.Item(22, 3).Value = 0.2738298704788917
$ws.Cells.Item(22, 4).Value = 0.02875303839744703
$ws.Cells.Item(22, 5).Value = 0.08607604712291028
$ws.Cells.Item(22, 6).Value = 0.7382852762447101
$ws.Cells.Item(22, 9).Value = 0.5670694740965772
$ws.Cells.Item(22, 12).Value = 0.279969982454233
$ws.Cells.Item(22, 14).Value = 1.047540108012008
$ws.Cells.Item(22, 15).Value = 2.410352654854051
$ws.Cells.Item(23, 2).Value = 1.594902577554763
$ws.Cells.Item(23, 3).Value = 0.2674935570931609
$ws.Cells.Item(23, 4).Value = 0.02813108915956519
$ws.Cells.Item(23, 5).Value = 0.08620102866204959
$ws.Cells.Item(23, 6).Value = 0.7334158779203079
$ws.Cells.Item(23, 9).Value = 0.5680509199981429
$ws.Cells.Item(23, 12).Value = 0.2748565084738459
$ws.Cells.Item(23, 14).Value = 1.049723701026494
$ws.Cells.Item(23, 15).Value = 2.401094851227867
$ws.Cells.Item(24, 2).Value = 1.405837344783777
$ws.Cells.Item(24, 3).Value = 0.243374003482387
$ws.Cells.Item(24, 4).Value = 0.02576354510038215
$ws.Cells.Item(24, 5).Value = 0.08675302885776226
$ws.Cells.Item(24, 6).Value = 0.7158595137828314
$ws.Cells.Item(24, 9).Value = 0.5724565881233445
$ws.Cells.Item(24, 12).Value = 0.2556855447520974
$ws.Cells.Item(24, 14).Value = 1.058766412722925
$ws.Cells.Item(24, 15).Value = 2.368983570557162
$ws.Cells.Item(25, 2).Value = 1.202147593333734
$ws.Cells.Item(25, 3).Value = 0.2171042145760111
$ws.Cells.Item(25, 4).Value = 0.02318490598914735
$ws.Cells.Item(25, 5).Value = 0.08751650656248344
$ws.Cells.Item(25, 6).Value = 0.6988089389161587
$ws.Cells.Item(25, 9).Value = 0.5786745710130496
$ws.Cells.Item(25, 12).Value = 0.2354273560700477
$ws.Cells.Item(25, 14).Value = 1.070173079762526
$ws.Cells.Item(25, 15).Value = 2.340629363084759

Write-Output "Updated 216 cells in pl_mw sheet (Case_5_20, 380 kV)"
